$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.730.57"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "1.875.82"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.87"
$ws.Range("E5").Value = "  +3.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("E7").Value = "  +6.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3963"
$ws.Range("E8").Value = "  +3.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.51"
$ws.Range("E9").Value = "  -2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08033"
$ws.Range("E10").Value = "  +2.72%  "

$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.89"
$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").Value = "1.882.88"
$ws.Range("E13").Value = "  +2.22%  "

$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.164"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("E17").Value = "  +3.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.20"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06629"
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "27.751.59"
$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.499"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.301"
$ws.Range("E25").Value = "  +1.76%  "

$ws.Range("D26").Value = "2.103.06"
$ws.Range("E26").Value = "  +1.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.61"
$ws.Range("E27").Value = "  +3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("E28").Value = "  +5.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").Value = "  +2.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.591"
$ws.Range("E30").Value = "  +2.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.55"
$ws.Range("E31").Value = "  +2.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9730"
$ws.Range("E32").Value = "  +4.98%  "

$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.456"
$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.633"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.313"
$ws.Range("E36").Value = "  +2.02%  "

$ws.Range("E37").Value = "  +2.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06114"
$ws.Range("E38").Value = "  +2.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.229"
$ws.Range("E39").Value = "  +2.03%  "

$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5997"
$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1907"
$ws.Range("E43").Value = "  +3.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.29"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5721"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.251"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.31"
$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.414"
$ws.Range("E48").Value = "  +1.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.937"
$ws.Range("E49").Value = "  +1.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06818"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.91"
$ws.Range("E51").Value = "  +4.59%  "
